$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert 4 new data rows (6-9) for the new "Regional Identity Service -
#    Australia East 3/4/5/6" entries, right after the existing "...East 2"
#    row (row 5). This pushes every following row down by 4.
# ---------------------------------------------------------------------------
$ws.Range("A6:A9").EntireRow.Insert()

# Clone the look (style/border) of the row just above into the freshly
# inserted blank rows - copy only columns A:B (column C is never used in
# this sheet and D gets its own formula below).
$ws.Range("A5:B5").Copy()
$ws.Range("A6:B9").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Fix up the two IP addresses that moved / changed for the existing
#    "Australia East 1" / "Australia East 2" regional identity rows.
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = "13.75.145.145"
$ws.Range("B5").Value = "40.82.217.103"

# ---------------------------------------------------------------------------
# 3. Populate the newly inserted rows with the new service entries.
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "Regional Identity Service - Australia East 3"
$ws.Range("B6").Value = "20.188.213.113"

$ws.Range("A7").Value = "Regional Identity Service - Australia East 4"
$ws.Range("B7").Value = "104.210.88.194"

$ws.Range("A8").Value = "Regional Identity Service - Australia East 5"
$ws.Range("B8").Value = "40.81.62.114"

$ws.Range("A9").Value = "Regional Identity Service - Australia East 6"
$ws.Range("B9").Value = "20.37.194.0/24"

# Give each new row the same "@{ serviceName = ...; IPs = ...; }," formula
# (as a shared formula, same family as the rest of column D) and restore
# the intended row height afterwards (writing a multi-line formula makes
# the engine auto-fit the row taller).
for ($r = 6; $r -le 9; $r++) {
    $formula = "=`"@{" + [char]10 + "    serviceName = '`"&A$r&`"'; " + [char]10 + "    IPs = '`"&B$r&`"';" + [char]10 + "},`""
    $ws.Range("D$r").Formula = $formula
}
for ($r = 6; $r -le 9; $r++) {
    $ws.Rows.Item($r).RowHeight = 14.65
}

# ---------------------------------------------------------------------------
# 4. Append two brand-new rows at the bottom of the table (after the last
#    existing row, "Azure Artifacts Blob - Australia East 3", which is now
#    row 27) for "Test Plans" and "Analytics service".
# ---------------------------------------------------------------------------
$ws.Range("A28:A29").EntireRow.Insert()

# These two new trailing rows use style 3 for BOTH the A and the B cell
# (unlike the rest of the table, where B uses style 4).
$ws.Range("A27").Copy()
$ws.Range("A28:B29").PasteSpecial(-4122)

$ws.Range("A28").Value = "Test Plans - Australia East 1"
$ws.Range("B28").Value = "20.40.177.101"

$ws.Range("A29").Value = "Analytics service - Australia East 1"
$ws.Range("B29").Value = "20.40.179.159"

for ($r = 28; $r -le 29; $r++) {
    $formula = "=`"@{" + [char]10 + "    serviceName = '`"&A$r&`"'; " + [char]10 + "    IPs = '`"&B$r&`"';" + [char]10 + "},`""
    $ws.Range("D$r").Formula = $formula
}
for ($r = 28; $r -le 29; $r++) {
    $ws.Rows.Item($r).RowHeight = 14.65
}

# ---------------------------------------------------------------------------
# 5. Selection now spans the full (bigger) table.
# ---------------------------------------------------------------------------
$ws.Range("A1:B29").Select()

Write-Host "Azure DevOps IP list updated"
